# Auto-generated script to update Leve market-price data cells
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets, per the scheduled
# market data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4828
$ws.Range("I62").Value = 3972.5
$ws.Range("J62").Value = 5968.6665
$ws.Range("K62").Value = 3972.5
$ws.Range("L62").Value = 5968.6665
$ws.Range("M62").Value = -3348.5
$ws.Range("N62").Value = -7216.6665
$ws.Range("H65").Value = 4828
$ws.Range("I65").Value = 3972.5
$ws.Range("J65").Value = 5968.6665
$ws.Range("K65").Value = 19862.5
$ws.Range("L65").Value = 29843.3325
$ws.Range("M65").Value = -16742.5
$ws.Range("N65").Value = -36083.3325
$ws.Range("H98").Value = 2405.0193
$ws.Range("I98").Value = 1757.8462
$ws.Range("J98").Value = 4346.5386
$ws.Range("K98").Value = 1757.8462
$ws.Range("L98").Value = 4346.5386
$ws.Range("M98").Value = -259.8462
$ws.Range("N98").Value = -7342.5386
$ws.Range("H112").Value = 1625.5555
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1633.7079
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 4901.1237
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -7117.1237
$ws.Range("H113").Value = 8597.308000000001
$ws.Range("I113").Value = 3973.2222
$ws.Range("J113").Value = 19001.5
$ws.Range("K113").Value = 3973.2222
$ws.Range("L113").Value = 19001.5
$ws.Range("M113").Value = -719.2222000000002
$ws.Range("N113").Value = -25509.5
$ws.Range("H122").Value = 2405.0193
$ws.Range("I122").Value = 1757.8462
$ws.Range("J122").Value = 4346.5386
$ws.Range("K122").Value = 5273.5386
$ws.Range("L122").Value = 13039.6158
$ws.Range("M122").Value = -2823.5386
$ws.Range("N122").Value = -17939.6158
$ws.Range("H123").Value = 42980
$ws.Range("J123").Value = 42980
$ws.Range("L123").Value = 42980
$ws.Range("N123").Value = -52780
$ws.Range("H129").Value = 897.119
$ws.Range("J129").Value = 964.43243
$ws.Range("L129").Value = 2893.29729
$ws.Range("N129").Value = -12893.29729
$ws.Range("H141").Value = 2119.946
$ws.Range("I141").Value = 1681.2667
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 5043.800099999999
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 136.1999000000005
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7649.8184
$ws.Range("I32").Value = 5560.3857
$ws.Range("J32").Value = 11491.678
$ws.Range("K32").Value = 5560.3857
$ws.Range("L32").Value = 11491.678
$ws.Range("M32").Value = -5273.3857
$ws.Range("N32").Value = -12065.678
$ws.Range("H137").Value = 51780
$ws.Range("J137").Value = 51780
$ws.Range("L137").Value = 51780
$ws.Range("N137").Value = -61980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1437.4242
$ws.Range("I20").Value = 1168.8422
$ws.Range("J20").Value = 1801.9286
$ws.Range("K20").Value = 1168.8422
$ws.Range("L20").Value = 1801.9286
$ws.Range("M20").Value = -921.8422
$ws.Range("N20").Value = -2295.9286
$ws.Range("H94").Value = 1836.8462
$ws.Range("I94").Value = 1930.75
$ws.Range("J94").Value = 710
$ws.Range("K94").Value = 1930.75
$ws.Range("L94").Value = 710
$ws.Range("M94").Value = -1479.75
$ws.Range("N94").Value = -1612
$ws.Range("H105").Value = 2838.0645
$ws.Range("I105").Value = 2841.0344
$ws.Range("K105").Value = 2841.0344
$ws.Range("M105").Value = -1094.0344
$ws.Range("H107").Value = 533.4231
$ws.Range("I107").Value = 543.4400000000001
$ws.Range("J107").Value = 283
$ws.Range("K107").Value = 543.4400000000001
$ws.Range("L107").Value = 283
$ws.Range("M107").Value = 1376.56
$ws.Range("N107").Value = -4123

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 8900
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 8900
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 26700
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -27512
$ws.Range("H85").Value = 8900
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 8900
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 26700
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -29508

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5165.263
$ws.Range("I70").Value = 5774.6665
$ws.Range("J70").Value = 4947.619
$ws.Range("K70").Value = 5774.6665
$ws.Range("L70").Value = 4947.619
$ws.Range("M70").Value = -5504.6665
$ws.Range("N70").Value = -5487.619
$ws.Range("H73").Value = 5165.263
$ws.Range("I73").Value = 5774.6665
$ws.Range("J73").Value = 4947.619
$ws.Range("K73").Value = 5774.6665
$ws.Range("L73").Value = 4947.619
$ws.Range("M73").Value = -4838.6665
$ws.Range("N73").Value = -6819.619

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5151.385
$ws.Range("I7").Value = 2194.6667
$ws.Range("J7").Value = 7685.7144
$ws.Range("K7").Value = 2194.6667
$ws.Range("L7").Value = 7685.7144
$ws.Range("M7").Value = -2082.6667
$ws.Range("N7").Value = -7909.7144
$ws.Range("H122").Value = 3419
$ws.Range("I122").Value = 3038.7646
$ws.Range("J122").Value = 4855.4443
$ws.Range("K122").Value = 9116.293799999999
$ws.Range("L122").Value = 14566.3329
$ws.Range("M122").Value = -6666.293799999999
$ws.Range("N122").Value = -19466.3329
$ws.Range("H126").Value = 5151.385
$ws.Range("I126").Value = 2194.6667
$ws.Range("J126").Value = 7685.7144
$ws.Range("K126").Value = 6584.000100000001
$ws.Range("L126").Value = 23057.1432
$ws.Range("M126").Value = -4114.000100000001
$ws.Range("N126").Value = -27997.1432
$ws.Range("H140").Value = 61552.867
$ws.Range("J140").Value = 61552.867
$ws.Range("L140").Value = 61552.867
$ws.Range("N140").Value = -71912.867

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 31321904
$ws.Range("I62").Value = 55558740
$ws.Range("K62").Value = 55558740
$ws.Range("M62").Value = -55558116
$ws.Range("H65").Value = 31321904
$ws.Range("I65").Value = 55558740
$ws.Range("K65").Value = 277793700
$ws.Range("M65").Value = -277790580
$ws.Range("H113").Value = 334.92856
$ws.Range("I113").Value = 283.76923
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 851.30769
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1318.69231
$ws.Range("N113").Value = -7340
$ws.Range("H136").Value = 4307.25
$ws.Range("I136").Value = 875
$ws.Range("K136").Value = 2625
$ws.Range("M136").Value = -75
